$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Changed Values")

$ws.Cells.Item(2, 1).Value = 39282
$ws.Cells.Item(2, 2).Value = 4
$ws.Cells.Item(2, 3).Value = 9
$ws.Cells.Item(2, 4).Value = "x1"
$ws.Cells.Item(3, 1).Value = 31550
$ws.Cells.Item(3, 2).Value = 5
$ws.Cells.Item(3, 3).Value = 10
$ws.Cells.Item(3, 4).Value = "x1"
$ws.Cells.Item(4, 1).Value = 32221
$ws.Cells.Item(4, 2).Value = 6
$ws.Cells.Item(4, 3).Value = 11
$ws.Cells.Item(4, 4).Value = "x1"
$ws.Cells.Item(5, 1).Value = 32673
$ws.Cells.Item(5, 2).Value = 7
$ws.Cells.Item(5, 3).Value = 12
$ws.Cells.Item(5, 4).Value = "x1"
$ws.Cells.Item(6, 1).Value = 33873
$ws.Cells.Item(6, 2).Value = 8
$ws.Cells.Item(6, 3).Value = 13
$ws.Cells.Item(6, 4).Value = "x1"
$ws.Cells.Item(7, 1).Value = 35545
$ws.Cells.Item(7, 2).Value = 9
$ws.Cells.Item(7, 3).Value = 14
$ws.Cells.Item(7, 4).Value = "x1"
$ws.Cells.Item(8, 1).Value = 36835
$ws.Cells.Item(8, 2).Value = 10
$ws.Cells.Item(8, 3).Value = 15
$ws.Cells.Item(8, 4).Value = "x1"
$ws.Cells.Item(9, 1).Value = 37525
$ws.Cells.Item(9, 2).Value = 11
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = "x1"
$ws.Cells.Item(10, 1).Value = 38382
$ws.Cells.Item(10, 2).Value = 12
$ws.Cells.Item(10, 3).Value = 17
$ws.Cells.Item(10, 4).Value = "x1"
$ws.Cells.Item(11, 1).Value = 40122
$ws.Cells.Item(11, 2).Value = 13
$ws.Cells.Item(11, 3).Value = 18
$ws.Cells.Item(11, 4).Value = "x1"
$ws.Cells.Item(12, 1).Value = 39818
$ws.Cells.Item(12, 2).Value = 14
$ws.Cells.Item(12, 3).Value = 19
$ws.Cells.Item(12, 4).Value = "x1"
$ws.Cells.Item(13, 1).Value = 149462
$ws.Cells.Item(13, 2).Value = 15
$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = "x1"
$ws.Cells.Item(14, 1).Value = 71062
$ws.Cells.Item(14, 2).Value = 16
$ws.Cells.Item(14, 3).Value = 8
$ws.Cells.Item(14, 4).Value = "x1"
$ws.Cells.Item(15, 1).Value = 31946
$ws.Cells.Item(15, 2).Value = 17
$ws.Cells.Item(15, 3).Value = 4
$ws.Cells.Item(15, 4).Value = "x1"
$ws.Cells.Item(16, 1).Value = 90208
$ws.Cells.Item(16, 2).Value = 18
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = "x1"
$ws.Cells.Item(17, 1).Value = 49082
$ws.Cells.Item(17, 2).Value = 19
$ws.Cells.Item(17, 3).Value = 6
$ws.Cells.Item(17, 4).Value = "x1"
$ws.Cells.Item(18, 1).Value = 1036991
$ws.Cells.Item(18, 2).Value = 20
$ws.Cells.Item(18, 3).Value = 31
$ws.Cells.Item(18, 4).Value = "x1"
$ws.Cells.Item(19, 1).Value = 64951
$ws.Cells.Item(19, 2).Value = 21
$ws.Cells.Item(19, 3).Value = 34
$ws.Cells.Item(19, 4).Value = "x1"
$ws.Cells.Item(20, 1).Value = 40620
$ws.Cells.Item(20, 2).Value = 22
$ws.Cells.Item(20, 3).Value = 20
$ws.Cells.Item(20, 4).Value = "x1"
$ws.Cells.Item(21, 1).Value = 33447
$ws.Cells.Item(21, 2).Value = 23
$ws.Cells.Item(21, 3).Value = 21
$ws.Cells.Item(21, 4).Value = "x1"
$ws.Cells.Item(22, 1).Value = 1265711
$ws.Cells.Item(22, 2).Value = 24
$ws.Cells.Item(22, 3).Value = 22
$ws.Cells.Item(22, 4).Value = "x1"
$ws.Cells.Item(23, 1).Value = 38235
$ws.Cells.Item(23, 2).Value = 25
$ws.Cells.Item(23, 3).Value = 23
$ws.Cells.Item(23, 4).Value = "x1"
$ws.Cells.Item(24, 1).Value = 38453
$ws.Cells.Item(24, 2).Value = 26
$ws.Cells.Item(24, 3).Value = 24
$ws.Cells.Item(24, 4).Value = "x1"
$ws.Cells.Item(25, 1).Value = 39392
$ws.Cells.Item(25, 2).Value = 27
$ws.Cells.Item(25, 3).Value = 25
$ws.Cells.Item(25, 4).Value = "x1"
$ws.Cells.Item(26, 1).Value = 39343
$ws.Cells.Item(26, 2).Value = 28
$ws.Cells.Item(26, 3).Value = 26
$ws.Cells.Item(26, 4).Value = "x1"
$ws.Cells.Item(27, 1).Value = 31665
$ws.Cells.Item(27, 2).Value = 29
$ws.Cells.Item(27, 3).Value = 27
$ws.Cells.Item(27, 4).Value = "x1"
$ws.Cells.Item(28, 1).Value = 39745
$ws.Cells.Item(28, 2).Value = 30
$ws.Cells.Item(28, 3).Value = 28
$ws.Cells.Item(28, 4).Value = "x1"
$ws.Cells.Item(29, 1).Value = 38611
$ws.Cells.Item(29, 2).Value = 31
$ws.Cells.Item(29, 3).Value = 29
$ws.Cells.Item(29, 4).Value = "x1"
$ws.Cells.Item(30, 1).Value = 42220
$ws.Cells.Item(30, 2).Value = 32
$ws.Cells.Item(30, 3).Value = 30
$ws.Cells.Item(30, 4).Value = "x1"
$ws.Cells.Item(31, 1).Value = 149448
$ws.Cells.Item(31, 2).Value = 33
$ws.Cells.Item(31, 3).Value = 32
$ws.Cells.Item(31, 4).Value = "x1"
$ws.Cells.Item(32, 1).Value = 32212
$ws.Cells.Item(32, 2).Value = 34
$ws.Cells.Item(32, 3).Value = 33
$ws.Cells.Item(32, 4).Value = "x1"
